$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "'20"
$ws.Range("C2").Style = "Normal"
$ws.Range("D2").Value = "'40000.00"
$ws.Range("D2").Style = "Normal"
$ws.Range("C3").Value = "'81"
$ws.Range("C3").Style = "Normal"
$ws.Range("D3").Value = "'288198.00"
$ws.Range("D3").Style = "Normal"
$ws.Range("C5").Value = "'129"
$ws.Range("C5").Style = "Normal"
$ws.Range("D5").Value = "'357386.40"
$ws.Range("D5").Style = "Normal"
$ws.Range("C6").Value = "'373"
$ws.Range("C6").Style = "Normal"
$ws.Range("D6").Value = "'995010.82"
$ws.Range("D6").Style = "Normal"
$ws.Range("C7").Value = "'67"
$ws.Range("C7").Style = "Normal"
$ws.Range("D7").Value = "'165289.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("C8").Value = "'741"
$ws.Range("C8").Style = "Normal"
$ws.Range("D8").Value = "'2710893.86"
$ws.Range("D8").Style = "Normal"
$ws.Range("C9").Value = "'24"
$ws.Range("C9").Style = "Normal"
$ws.Range("D9").Value = "'65100.00"
$ws.Range("D9").Style = "Normal"
$ws.Range("C11").Value = "'34"
$ws.Range("C11").Style = "Normal"
$ws.Range("D11").Value = "'86177.00"
$ws.Range("D11").Style = "Normal"
$ws.Range("C12").Value = "'148"
$ws.Range("C12").Style = "Normal"
$ws.Range("D12").Value = "'433575.18"
$ws.Range("D12").Style = "Normal"
$ws.Range("C13").Value = "'77"
$ws.Range("C13").Style = "Normal"
$ws.Range("D13").Value = "'193800.00"
$ws.Range("D13").Style = "Normal"
$ws.Range("C14").Value = "'91"
$ws.Range("C14").Style = "Normal"
$ws.Range("D14").Value = "'223788.98"
$ws.Range("D14").Style = "Normal"
$ws.Range("C15").Value = "'18"
$ws.Range("C15").Style = "Normal"
$ws.Range("D15").Value = "'39593.58"
$ws.Range("D15").Style = "Normal"
$ws.Range("C16").Value = "'120"
$ws.Range("C16").Style = "Normal"
$ws.Range("D16").Value = "'537507.26"
$ws.Range("D16").Style = "Normal"
$ws.Range("C17").Value = "'173"
$ws.Range("C17").Style = "Normal"
$ws.Range("D17").Value = "'396444.87"
$ws.Range("D17").Style = "Normal"
$ws.Range("C18").Value = "'11"
$ws.Range("C18").Style = "Normal"
$ws.Range("D18").Value = "'30500.00"
$ws.Range("D18").Style = "Normal"
$ws.Range("C23").Value = "'297"
$ws.Range("C23").Style = "Normal"
$ws.Range("D23").Value = "'1158826.10"
$ws.Range("D23").Style = "Normal"
$ws.Range("C26").Value = "'38"
$ws.Range("C26").Style = "Normal"
$ws.Range("D26").Value = "'114572.00"
$ws.Range("D26").Style = "Normal"
$ws.Range("C30").Value = "'42"
$ws.Range("C30").Style = "Normal"
$ws.Range("D30").Value = "'187060.00"
$ws.Range("D30").Style = "Normal"
$ws.Range("C31").Value = "'52"
$ws.Range("C31").Style = "Normal"
$ws.Range("D31").Value = "'123310.00"
$ws.Range("D31").Style = "Normal"
$ws.Range("C37").Value = "'371"
$ws.Range("C37").Style = "Normal"
$ws.Range("D37").Value = "'1456190.10"
$ws.Range("D37").Style = "Normal"
$ws.Range("C45").Value = "'46"
$ws.Range("C45").Style = "Normal"
$ws.Range("D45").Value = "'157715.92"
$ws.Range("D45").Style = "Normal"
$ws.Range("C50").Value = "'113"
$ws.Range("C50").Style = "Normal"
$ws.Range("D50").Value = "'288178.33"
$ws.Range("D50").Style = "Normal"
$ws.Range("C74").Value = "'19"
$ws.Range("C74").Style = "Normal"
$ws.Range("D74").Value = "'77500.00"
$ws.Range("D74").Style = "Normal"
$ws.Range("C75").Value = "'41"
$ws.Range("C75").Style = "Normal"
$ws.Range("D75").Value = "'116579.25"
$ws.Range("D75").Style = "Normal"
$ws.Range("C77").Value = "'88"
$ws.Range("C77").Style = "Normal"
$ws.Range("D77").Value = "'232987.00"
$ws.Range("D77").Style = "Normal"
$ws.Range("C78").Value = "'204"
$ws.Range("C78").Style = "Normal"
$ws.Range("D78").Value = "'565693.00"
$ws.Range("D78").Style = "Normal"
$ws.Range("C80").Value = "'476"
$ws.Range("C80").Style = "Normal"
$ws.Range("D80").Value = "'2044100.03"
$ws.Range("D80").Style = "Normal"
$ws.Range("C84").Value = "'70"
$ws.Range("C84").Style = "Normal"
$ws.Range("D84").Value = "'236657.55"
$ws.Range("D84").Style = "Normal"
$ws.Range("C85").Value = "'34"
$ws.Range("C85").Style = "Normal"
$ws.Range("D85").Value = "'106669.00"
$ws.Range("D85").Style = "Normal"
$ws.Range("C88").Value = "'68"
$ws.Range("C88").Style = "Normal"
$ws.Range("D88").Value = "'300956.08"
$ws.Range("D88").Style = "Normal"
$ws.Range("C89").Value = "'106"
$ws.Range("C89").Style = "Normal"
$ws.Range("D89").Value = "'274120.00"
$ws.Range("D89").Style = "Normal"
$ws.Range("C90").Value = "'38"
$ws.Range("C90").Style = "Normal"
$ws.Range("D90").Value = "'79555.00"
$ws.Range("D90").Style = "Normal"
$ws.Range("C91").Value = "'56"
$ws.Range("C91").Style = "Normal"
$ws.Range("D91").Value = "'160957.14"
$ws.Range("D91").Style = "Normal"
$ws.Range("C92").Value = "'63"
$ws.Range("C92").Style = "Normal"
$ws.Range("D92").Value = "'168445.71"
$ws.Range("D92").Style = "Normal"
$ws.Range("C93").Value = "'119"
$ws.Range("C93").Style = "Normal"
$ws.Range("D93").Value = "'306355.00"
$ws.Range("D93").Style = "Normal"
$ws.Range("C94").Value = "'23"
$ws.Range("C94").Style = "Normal"
$ws.Range("D94").Value = "'47500.00"
$ws.Range("D94").Style = "Normal"
$ws.Range("C95").Value = "'138"
$ws.Range("C95").Style = "Normal"
$ws.Range("D95").Value = "'399486.00"
$ws.Range("D95").Style = "Normal"
$ws.Range("C97").Value = "'9"
$ws.Range("C97").Style = "Normal"
$ws.Range("D97").Value = "'19500.00"
$ws.Range("D97").Style = "Normal"
$ws.Range("C99").Value = "'47"
$ws.Range("C99").Style = "Normal"
$ws.Range("D99").Value = "'138000.00"
$ws.Range("D99").Style = "Normal"
$ws.Range("C100").Value = "'65"
$ws.Range("C100").Style = "Normal"
$ws.Range("D100").Value = "'157479.00"
$ws.Range("D100").Style = "Normal"
$ws.Range("C103").Value = "'23"
$ws.Range("C103").Style = "Normal"
$ws.Range("D103").Value = "'59830.00"
$ws.Range("D103").Style = "Normal"
$ws.Range("C104").Value = "'53"
$ws.Range("C104").Style = "Normal"
$ws.Range("D104").Value = "'107500.00"
$ws.Range("D104").Style = "Normal"
$ws.Range("C122").Value = "'242"
$ws.Range("C122").Style = "Normal"
$ws.Range("D122").Value = "'667508.00"
$ws.Range("D122").Style = "Normal"
$ws.Range("C123").Value = "'106"
$ws.Range("C123").Style = "Normal"
$ws.Range("D123").Value = "'280081.45"
$ws.Range("D123").Style = "Normal"
$ws.Range("C124").Value = "'469"
$ws.Range("C124").Style = "Normal"
$ws.Range("D124").Value = "'2062632.06"
$ws.Range("D124").Style = "Normal"
$ws.Range("C129").Value = "'42"
$ws.Range("C129").Style = "Normal"
$ws.Range("D129").Value = "'156579.76"
$ws.Range("D129").Style = "Normal"
$ws.Range("C132").Value = "'84"
$ws.Range("C132").Style = "Normal"
$ws.Range("D132").Value = "'379163.75"
$ws.Range("D132").Style = "Normal"
$ws.Range("C135").Value = "'215"
$ws.Range("C135").Style = "Normal"
$ws.Range("D135").Value = "'609550.00"
$ws.Range("D135").Style = "Normal"
$ws.Range("C138").Value = "'573"
$ws.Range("C138").Style = "Normal"
$ws.Range("D138").Value = "'1442546.00"
$ws.Range("D138").Style = "Normal"
$ws.Range("C139").Value = "'1827"
$ws.Range("C139").Style = "Normal"
$ws.Range("D139").Value = "'4932715.93"
$ws.Range("D139").Style = "Normal"
$ws.Range("C140").Value = "'2687"
$ws.Range("C140").Style = "Normal"
$ws.Range("D140").Value = "'6811637.46"
$ws.Range("D140").Style = "Normal"
$ws.Range("C141").Value = "'2673"
$ws.Range("C141").Style = "Normal"
$ws.Range("D141").Value = "'11619111.82"
$ws.Range("D141").Style = "Normal"
$ws.Range("C142").Value = "'360"
$ws.Range("C142").Style = "Normal"
$ws.Range("D142").Value = "'1023872.94"
$ws.Range("D142").Style = "Normal"
$ws.Range("C143").Value = "'128"
$ws.Range("C143").Style = "Normal"
$ws.Range("D143").Value = "'321984.00"
$ws.Range("D143").Style = "Normal"
$ws.Range("C144").Value = "'256"
$ws.Range("C144").Style = "Normal"
$ws.Range("D144").Value = "'681516.33"
$ws.Range("D144").Style = "Normal"
$ws.Range("C145").Value = "'1053"
$ws.Range("C145").Style = "Normal"
$ws.Range("D145").Value = "'2776892.25"
$ws.Range("D145").Style = "Normal"
$ws.Range("C146").Value = "'511"
$ws.Range("C146").Style = "Normal"
$ws.Range("D146").Value = "'1518291.49"
$ws.Range("D146").Style = "Normal"
$ws.Range("C147").Value = "'379"
$ws.Range("C147").Style = "Normal"
$ws.Range("D147").Value = "'957377.83"
$ws.Range("D147").Style = "Normal"
$ws.Range("C148").Value = "'148"
$ws.Range("C148").Style = "Normal"
$ws.Range("D148").Value = "'367500.00"
$ws.Range("D148").Style = "Normal"
$ws.Range("C149").Value = "'429"
$ws.Range("C149").Style = "Normal"
$ws.Range("D149").Value = "'1398905.46"
$ws.Range("D149").Style = "Normal"
$ws.Range("C150").Value = "'860"
$ws.Range("C150").Style = "Normal"
$ws.Range("D150").Value = "'2094695.82"
$ws.Range("D150").Style = "Normal"
$ws.Range("C194").Value = "'55"
$ws.Range("C194").Style = "Normal"
$ws.Range("D194").Value = "'176300.00"
$ws.Range("D194").Style = "Normal"
$ws.Range("C197").Value = "'352"
$ws.Range("C197").Style = "Normal"
$ws.Range("D197").Value = "'949788.00"
$ws.Range("D197").Style = "Normal"
$ws.Range("C199").Value = "'649"
$ws.Range("C199").Style = "Normal"
$ws.Range("D199").Value = "'2443424.58"
$ws.Range("D199").Style = "Normal"
$ws.Range("C203").Value = "'154"
$ws.Range("C203").Style = "Normal"
$ws.Range("D203").Value = "'469133.00"
$ws.Range("D203").Style = "Normal"
$ws.Range("C204").Value = "'63"
$ws.Range("C204").Style = "Normal"
$ws.Range("D204").Value = "'183926.00"
$ws.Range("D204").Style = "Normal"
$ws.Range("C205").Value = "'80"
$ws.Range("C205").Style = "Normal"
$ws.Range("D205").Value = "'190720.00"
$ws.Range("D205").Style = "Normal"
$ws.Range("C207").Value = "'122"
$ws.Range("C207").Style = "Normal"
$ws.Range("D207").Value = "'570888.14"
$ws.Range("D207").Style = "Normal"
$ws.Range("C240").Value = "'82"
$ws.Range("C240").Style = "Normal"
$ws.Range("D240").Value = "'228538.00"
$ws.Range("D240").Style = "Normal"
$ws.Range("C241").Value = "'145"
$ws.Range("C241").Style = "Normal"
$ws.Range("D241").Value = "'374200.00"
$ws.Range("D241").Style = "Normal"
$ws.Range("C242").Value = "'494"
$ws.Range("C242").Style = "Normal"
$ws.Range("D242").Value = "'1281575.83"
$ws.Range("D242").Style = "Normal"
$ws.Range("C243").Value = "'85"
$ws.Range("C243").Style = "Normal"
$ws.Range("D243").Value = "'246627.11"
$ws.Range("D243").Style = "Normal"
$ws.Range("C244").Value = "'959"
$ws.Range("C244").Style = "Normal"
$ws.Range("D244").Value = "'3447071.19"
$ws.Range("D244").Style = "Normal"
$ws.Range("C247").Value = "'79"
$ws.Range("C247").Style = "Normal"
$ws.Range("D247").Value = "'180500.00"
$ws.Range("D247").Style = "Normal"
$ws.Range("C252").Value = "'120"
$ws.Range("C252").Style = "Normal"
$ws.Range("D252").Value = "'428812.14"
$ws.Range("D252").Style = "Normal"
$ws.Range("C253").Value = "'210"
$ws.Range("C253").Style = "Normal"
$ws.Range("D253").Value = "'461363.00"
$ws.Range("D253").Style = "Normal"
